$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.892.85"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.950.71"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.16"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.32"
$ws.Range("E6").Value = "  +9.50%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.93%  "
$ws.Range("D9").Value = "2.945.93"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.92"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "3.427.72"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.97"
$ws.Range("E17").Value = "  +9.60%  "
$ws.Range("D18").Value = "2.943.65"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "57.848.86"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.44"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  +5.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  +7.99%  "
$ws.Range("E23").Value = "  +7.77%  "
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.06"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.50"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +5.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.52"
$ws.Range("E30").Value = "  +6.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.48"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("E34").Value = "  +6.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.948"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.01"
$ws.Range("E37").Value = "  +8.16%  "
$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  +13.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.35"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +15.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "382.04"
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0348"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").Value = "2.706.86"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.29"
$ws.Range("E46").Value = "  +5.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.238"
$ws.Range("E47").Value = "  +4.65%  "
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.98"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.96"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("E51").Value = "  +4.13%  "
